# Auto update Excel log — append new sensor readings to the PIR, Humidity,
# and Temperature sheets (rows captured 2026-01-28 13:03:xx for the
# Bathroom sensors), matching the source system's ongoing data export.

$wb = $excel.ActiveWorkbook

# --- Data to append -------------------------------------------------------
# Each row: RowNumber, Date, Timestamp, Hour, Location, Value, Status

$pirRows = @(
    @(418, "2026-01-28", "13:03:02", "13:00", "Bathroom", "No Motion", "Inactive"),
    @(419, "2026-01-28", "13:03:05", "13:00", "Bathroom", "No Motion", "Inactive"),
    @(420, "2026-01-28", "13:03:09", "13:00", "Bathroom", "No Motion", "Inactive"),
    @(421, "2026-01-28", "13:03:14", "13:00", "Bathroom", "No Motion", "Inactive"),
    @(422, "2026-01-28", "13:03:19", "13:00", "Bathroom", "No Motion", "Inactive"),
    @(423, "2026-01-28", "13:03:26", "13:00", "Bathroom", "No Motion", "Inactive"),
    @(424, "2026-01-28", "13:03:30", "13:00", "Bathroom", "No Motion", "Inactive"),
    @(425, "2026-01-28", "13:03:34", "13:00", "Bathroom", "No Motion", "Inactive"),
    @(426, "2026-01-28", "13:03:39", "13:00", "Bathroom", "No Motion", "Inactive"),
    @(427, "2026-01-28", "13:03:46", "13:00", "Bathroom", "No Motion", "Inactive"),
    @(428, "2026-01-28", "13:03:50", "13:00", "Bathroom", "No Motion", "Inactive"),
    @(429, "2026-01-28", "13:03:54", "13:00", "Bathroom", "No Motion", "Inactive"),
    @(430, "2026-01-28", "13:03:59", "13:00", "Bathroom", "No Motion", "Inactive")
)

$humidityRows = @(
    @(391, "2026-01-28", "13:03:00", "13:00", "Bathroom", "88.4%", "Active"),
    @(392, "2026-01-28", "13:03:03", "13:00", "Bathroom", "88.3%", "Active"),
    @(393, "2026-01-28", "13:03:07", "13:00", "Bathroom", "88.4%", "Active"),
    @(394, "2026-01-28", "13:03:15", "13:00", "Bathroom", "88.3%", "Active"),
    @(395, "2026-01-28", "13:03:23", "13:00", "Bathroom", "87.5%", "Active"),
    @(396, "2026-01-28", "13:03:27", "13:00", "Bathroom", "88.4%", "Active"),
    @(397, "2026-01-28", "13:03:31", "13:00", "Bathroom", "87.5%", "Active"),
    @(398, "2026-01-28", "13:03:35", "13:00", "Bathroom", "88.4%", "Active"),
    @(399, "2026-01-28", "13:03:44", "13:00", "Bathroom", "88.5%", "Active"),
    @(400, "2026-01-28", "13:03:47", "13:00", "Bathroom", "88.5%", "Active"),
    @(401, "2026-01-28", "13:03:52", "13:00", "Bathroom", "87.5%", "Active")
)

$temperatureRows = @(
    @(391, "2026-01-28", "13:03:01", "13:00", "Bathroom", "22.8C", "Active"),
    @(392, "2026-01-28", "13:03:04", "13:00", "Bathroom", "22.7C", "Active"),
    @(393, "2026-01-28", "13:03:08", "13:00", "Bathroom", "22.8C", "Active"),
    @(394, "2026-01-28", "13:03:16", "13:00", "Bathroom", "22.7C", "Active"),
    @(395, "2026-01-28", "13:03:24", "13:00", "Bathroom", "22.7C", "Active"),
    @(396, "2026-01-28", "13:03:28", "13:00", "Bathroom", "22.8C", "Active"),
    @(397, "2026-01-28", "13:03:33", "13:00", "Bathroom", "22.8C", "Active"),
    @(398, "2026-01-28", "13:03:37", "13:00", "Bathroom", "22.8C", "Active"),
    @(399, "2026-01-28", "13:03:45", "13:00", "Bathroom", "22.8C", "Active"),
    @(400, "2026-01-28", "13:03:49", "13:00", "Bathroom", "22.8C", "Active"),
    @(401, "2026-01-28", "13:03:53", "13:00", "Bathroom", "22.7C", "Active")
)

function Append-SensorRows($SheetName, $Rows) {
    $ws = $wb.Worksheets.Item($SheetName)

    foreach ($row in $Rows) {
        $r = $row[0]

        # Columns A-C hold date/time-looking text, and column E sometimes
        # holds percentage-looking text (e.g. "88.4%") — all must stay
        # literal text, not auto-converted to Excel date/time/number values.
        $ws.Range("A$r").NumberFormat = "@"
        $ws.Range("B$r").NumberFormat = "@"
        $ws.Range("C$r").NumberFormat = "@"
        $ws.Range("E$r").NumberFormat = "@"

        $ws.Range("A$r").Value = $row[1]
        $ws.Range("B$r").Value = $row[2]
        $ws.Range("C$r").Value = $row[3]
        $ws.Range("D$r").Value = $row[4]
        $ws.Range("E$r").Value = $row[5]
        $ws.Range("F$r").Value = $row[6]
    }
}

Append-SensorRows "PIR" $pirRows
Append-SensorRows "Humidity" $humidityRows
Append-SensorRows "Temperature" $temperatureRows
